# Actualización automática del mapa (2025-07-21 07:31:02)
# Elimina la fila del caso 5700 (ESTOMBA 2119) que ya no corresponde,
# desplazando el resto de los registros hacia arriba.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25 contains case 5700 / ESTOMBA 2119 which must be removed entirely.
$ws.Rows.Item(25).Delete()
